$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Intro paragraph: merge the tail of the paragraph (", nos vamos a
#    adentrar ... se explicarán ... efectiva.") into a single run with new
#    wording ("se explicará ... se abordarán ... efectiva.").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
  ", nos vamos a adentrar en el lenguaje principal utilizado en React: JavaScript. A través de ejemplos prácticos, se explicarán los fundamentos de JavaScript necesarios para comprender y trabajar con React de manera efectiva.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  " se explicará el lenguaje principal utilizado en React: JavaScript. A través de ejemplos prácticos, se abordarán los fundamentos de JavaScript necesarios para comprender y trabajar con React de manera efectiva.",
  2) | Out-Null

# ---------------------------------------------------------------------------
# Helper paragraphs: rebuild four bullet paragraphs so they gain
# <w:proofErr> spellStart/spellEnd markers around "map"/"find"/"filter",
# "Fetch" and "APIs" (these also split the "Fetch API" text into separate
# runs, one of which keeps just "Fetch").
# ---------------------------------------------------------------------------

$pkgHeader = "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'><w:body>"
$pkgFooter = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# --- 2) "Exploración de funciones avanzadas como map, find, filter, etc." ---
$rng = $d.Content
$rng.Find.Execute("Exploración de funciones avanzadas como map, find, filter, etc.",
  $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null

$body = '<w:p w14:paraId="79E305EA" w14:textId="77777777" w:rsidR="00CD09D7" w:rsidRDefault="00CD09D7" w:rsidP="001F703B">' + `
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="13"/></w:numPr><w:contextualSpacing w:val="0"/><w:jc w:val="both"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Exploración de funciones avanzadas como </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F86A5C"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>map</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F86A5C"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>find</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F86A5C"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>filter</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>, etc.</w:t></w:r></w:p>'
$rng.InsertXML($pkgHeader + $body + $pkgFooter)

# --- 3) "Fetch API" heading: split into "Fetch" + " API" runs -------------
$rng = $d.Content
$rng.Find.Execute("Fetch API", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null

$body = '<w:p w14:paraId="197A6280" w14:textId="0AC714AE" w:rsidR="00CD09D7" w:rsidRPr="00220A8C" w:rsidRDefault="00CD09D7" w:rsidP="001F703B">' + `
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr><w:ind w:left="714" w:hanging="357"/><w:contextualSpacing w:val="0"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00220A8C"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Fetch</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r w:rsidRPr="00220A8C"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> API</w:t></w:r></w:p>'
$rng.InsertXML($pkgHeader + $body + $pkgFooter)

# --- 4) "Descripción de la Fetch API para realizar..." --------------------
$rng = $d.Content
$rng.Find.Execute("Descripción de la Fetch API para realizar peticiones HTTP desde el navegador.",
  $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null

$body = '<w:p w14:paraId="07219814" w14:textId="77777777" w:rsidR="00CD09D7" w:rsidRDefault="00CD09D7" w:rsidP="001F703B">' + `
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="13"/></w:numPr><w:contextualSpacing w:val="0"/><w:jc w:val="both"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Descripción de la </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Fetch</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> API para realizar peticiones HTTP desde el navegador.</w:t></w:r></w:p>'
$rng.InsertXML($pkgHeader + $body + $pkgFooter)

# --- 5) "Aplicación del Fetch API para consumir APIs y realizar..." -------
$rng = $d.Content
$rng.Find.Execute("Aplicación del Fetch API para consumir APIs y realizar peticiones al servidor backend.",
  $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null

$body = '<w:p w14:paraId="0489F095" w14:textId="1251932A" w:rsidR="00CD09D7" w:rsidRDefault="00CD09D7" w:rsidP="001F703B">' + `
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="13"/></w:numPr><w:contextualSpacing w:val="0"/><w:jc w:val="both"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Aplicación del </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Fetch</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> API para consumir </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>APIs</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> y realizar peticiones al servidor backend.</w:t></w:r></w:p>'
$rng.InsertXML($pkgHeader + $body + $pkgFooter)
